$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: text update (Datenbank optimiert -> Datenbank optimiert(Probleme mit einigen Tabellen und Relationen))
# and row grows taller to fit the extra wrapped text.
$ws.Rows(11).RowHeight = 54
$ws.Cells.Item(11, 5).Value = "Check-Klasse programmiert, Verification von einigen Methoden und programmierten Seiten, Flipchart und Burndownchart aktualisiert, Datenbank optimiert(Probleme mit einigen Tabellen und Relationen), Programmierbeihilfe"

# Row 12 keeps its original content (10.12.2019 / Sprintplanung beendet...) - no change needed.

# Row 13: previously blank, now a real entry (17.12.2019)
$ws.Rows(13).RowHeight = 27
$ws.Cells.Item(13, 1).Value = "Di"
$ws.Cells.Item(13, 2).Value2 = "17.12.2019"
$ws.Cells.Item(13, 3).Value = 0.32291666666666669
$ws.Cells.Item(13, 4).Value = 0.54861111111111105
$ws.Cells.Item(13, 5).Value = "Programmieren an Story Kursauflösung bzw. Löschen und Programmierbeihilfe, Testung von Seiten"

# Row 14: previously blank, now a real entry (07.01.2020)
$ws.Rows(14).RowHeight = 27
$ws.Cells.Item(14, 1).Value = "Di"
$ws.Cells.Item(14, 2).Value2 = "07.01.2020"
$ws.Cells.Item(14, 3).Value = 0.32291666666666669
$ws.Cells.Item(14, 4).Value = 0.54861111111111105
$ws.Cells.Item(14, 5).Value = "Story Kursauflösung beendet, Für Sprintfertigstellung gesorgt (Protokolle auf IST-Stand, Burndownchart etc.) "
